$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Craft")

# Row 59: Chest_block_item crafted from Wood Planks (was empty -> fixes furnace dupe bug)
$ws.Range("A59").Value = "Chest_block_item"
$ws.Range("B59").Value = 1
$ws.Range("C59").Value = 3
$ws.Range("D59").Value = 3
$ws.Range("E59").Value = "Wood_Planks_block_item,Wood_Planks_block_item,Wood_Planks_block_item"
$ws.Range("F59").Value = "Wood_Planks_block_item,null,Wood_Planks_block_item"
$ws.Range("G59").Value = "Wood_Planks_block_item,Wood_Planks_block_item,Wood_Planks_block_item"

# Row 60: Furnace_block_item crafted from Cobblestone (was empty -> fixes furnace dupe bug)
$ws.Range("A60").Value = "Furnace_block_item"
$ws.Range("B60").Value = 1
$ws.Range("C60").Value = 3
$ws.Range("D60").Value = 3
$ws.Range("E60").Value = "Cobblestone_block_item,Cobblestone_block_item,Cobblestone_block_item"
$ws.Range("F60").Value = "Cobblestone_block_item,null,Cobblestone_block_item"
$ws.Range("G60").Value = "Cobblestone_block_item,Cobblestone_block_item,Cobblestone_block_item"

# Row 61: Grass_block_item recipe variant (Dirt + Leaves)
$ws.Range("A61").Value = "Grass_block_item"
$ws.Range("B61").Value = 1
$ws.Range("C61").Value = 2
$ws.Range("D61").Value = 1
$ws.Range("E61").Value = "Dirt_block_item,Leaves_block_item"

# Row 62: Grass_block_item recipe variant (Leaves + Dirt)
$ws.Range("A62").Value = "Grass_block_item"
$ws.Range("B62").Value = 1
$ws.Range("C62").Value = 2
$ws.Range("D62").Value = 1
$ws.Range("E62").Value = "Leaves_block_item,Dirt_block_item"

# Row 63: new blank styled cell at G63 (center + underline), matching row layout shift
$ws.Range("G63").HorizontalAlignment = -4108
$ws.Range("G63").Font.Underline = 2

# Restore view/selection state
$wsItems = $wb.Worksheets.Item("Items")
$wsItems.Activate()
$wsItems.Range("A58").Select()

$ws.Activate()
$ws.Range("E63").Select()
